$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.720.92'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '1.700.02'
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = "'315.97"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').Value = "'1.004"
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Value = "'0.3935"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('D8').Value = "'0.4043"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.67%  '
$ws.Range('D9').Value = "'1.524"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('D10').Value = "'53.90"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.67%  '
$ws.Range('D11').Value = "'1.004"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').Value = "'0.08872"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('D13').Value = "'7.402"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.82%  '
$ws.Range('D14').Value = "'23.64"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').Value = "'8.192"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.84%  '
$ws.Range('D16').Value = "'0.00001324"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').Value = '1.711.56'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').Value = "'99.71"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.22%  '
$ws.Range('D19').Value = "'0.07056"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.59%  '
$ws.Range('D20').Value = "'19.70"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = "'7.085"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.84%  '
$ws.Range('D22').Value = "'1.007"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.58%  '
$ws.Range('E23').Value = '  +4.77%  '
$ws.Range('D24').Value = '24.722.64'
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('D25').Value = "'3.133"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.96%  '
$ws.Range('D26').Value = "'2.368"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('D27').Value = "'22.76"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.87%  '
$ws.Range('D28').Value = "'163.16"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.09%  '
$ws.Range('D29').Value = "'8.775"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +17.38%  '
$ws.Range('D30').Value = "'135.87"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.23%  '
$ws.Range('E31').Value = '  -0.67%  '
$ws.Range('D32').Value = "'0.09038"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.22%  '
$ws.Range('D33').Value = "'7.673"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.64%  '
$ws.Range('D34').Value = "'1.067"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.06%  '
$ws.Range('D35').Value = "'1.976"
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Value = "'11.14"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.27%  '
$ws.Range('D37').Value = "'0.2756"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('D38').Value = "'14.53"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.24%  '
$ws.Range('D39').Value = "'0.02779"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('D40').Value = "'0.09156"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.60%  '
$ws.Range('D41').Value = "'1.462"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').Value = "'0.7668"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('D43').Value = "'15.96"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.79%  '
$ws.Range('D44').Value = "'0.7178"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('D45').Value = "'2.573"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.59%  '
$ws.Range('D46').Value = "'4.217"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('D48').Value = "'1.337"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('D49').Value = "'140.14"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('D50').Value = "'90.98"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.07%  '
$ws.Range('D51').Value = "'0.07987"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.65%  '
